# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates column G ("K") values for rows 2-39 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 1
    6  = 1
    7  = 1
    8  = 0
    9  = 3
    10 = 1
    11 = 1
    12 = 2
    13 = 1
    14 = 0
    15 = 3
    16 = 0
    17 = 0
    18 = 1
    19 = 0
    21 = 0
    22 = 0
    23 = 0
    24 = 1
    25 = 1
    26 = 1
    27 = 1
    28 = 1
    29 = 1
    30 = 1
    31 = 2
    32 = 2
    33 = 0
    34 = 2
    35 = 1
    36 = 1
    37 = 1
    38 = 2
    39 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
